$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Add an 8th (blank) page: after the last paragraph of the body
#    (the tiny sz=12 spacer paragraph right before the final sectPr),
#    insert a page-break paragraph followed by a new trailing blank
#    paragraph, both sharing that same tiny formatting.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last

# Mint a brand-new paragraph after the existing last one. Word copies
# the paragraph formatting (pPr) from the source paragraph onto it.
$lastPara.Range.InsertParagraphAfter()

$breakPara = $d.Paragraphs.Last.Previous

# Fill that new paragraph with a run containing an explicit page break,
# keeping the same tiny spacer formatting on both pPr/rPr and the run.
$pageBreakXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="12"/><w:szCs w:val="12"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="12"/><w:szCs w:val="12"/><w:lang w:val="en-US"/></w:rPr><w:br w:type="page"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$breakPara.Range.InsertXML($pageBreakXml)

# ------------------------------------------------------------------
# 2) The cached result of the footer's PAGE field was last computed
#    as "4"; refresh it to "7" to match the now-longer document.
# ------------------------------------------------------------------
$footerRange = $d.Sections.Item(1).Footers.Item(1).Range
$footerRange.Find.Execute("4", $true, $false, $false, $false, $false, $true, 1, $false, "7", 2)
